$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D stays formatted as text so numeric-looking strings
# (e.g. "0.0081") are not auto-converted to numbers by Excel.
$ws.Range("D2:D11").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Move to location (11, 8) and remove the toolkit."
$ws.Range("B2").Value = 24.588919
$ws.Range("C2").Value = 3868
$ws.Range("D2").Value = "0.0081"
$ws.Range("E2").Value = "bf987438-11fe-4ea4-8ac8-aef65b60f040"

# Row 3
$ws.Range("A3").Value = "Move to location (7, 5) and remove the liquid spill."
$ws.Range("B3").Value = 25.774564
$ws.Range("C3").Value = 3844
$ws.Range("D3").Value = "0.00804"
$ws.Range("E3").Value = "43b8ad0a-efe5-4ecb-87e8-7a597de1268d"

# Row 4
$ws.Range("A4").Value = "Move to location (8, 6) and remove the large debris."
$ws.Range("B4").Value = 27.574987
$ws.Range("C4").Value = 3952
$ws.Range("D4").Value = "0.00855"
$ws.Range("E4").Value = "9e446b93-6741-4fdf-89e3-672158208eb6"

# Row 5
$ws.Range("A5").Value = "Move to location (2, 4) and remove the dust."
$ws.Range("B5").Value = 27.627885
$ws.Range("C5").Value = 3829
$ws.Range("D5").Value = "0.00813"
$ws.Range("E5").Value = "e9f4fe67-4b74-4bcb-9973-71b8787c0c5e"

# Row 6
$ws.Range("A6").Value = "Move to location (5, 2) and remove the grass."
$ws.Range("B6").Value = 23.637939
$ws.Range("C6").Value = 3839
$ws.Range("D6").Value = "0.00774"
$ws.Range("E6").Value = "f153fe00-ec80-4937-bae3-a58698a04708"

# Row 7
$ws.Range("A7").Value = "Move to location (6, 7) and remove the small debris."
$ws.Range("B7").Value = 28.864845
$ws.Range("C7").Value = 3960
$ws.Range("D7").Value = "0.00864"
$ws.Range("E7").Value = "2cad060a-57b4-4d08-9626-6a4ab989612a"

# Row 8
$ws.Range("A8").Value = "Move to location (3, 6) and remove the vehicle."
$ws.Range("B8").Value = 24.688907
$ws.Range("C8").Value = 3858
$ws.Range("D8").Value = "0.00855"
$ws.Range("E8").Value = "b118f4fd-a574-4486-afe5-d83b9ebf4959"

# Row 9
$ws.Range("A9").Value = "Move to location (6, 6) and remove the construction materials."
$ws.Range("B9").Value = 25.909734
$ws.Range("C9").Value = 3927
$ws.Range("D9").Value = "0.0084"
$ws.Range("E9").Value = "375dad59-fd8a-4822-b932-a14628afbed2"

# Row 10
$ws.Range("A10").Value = "Move to location (3, 9) and remove the tree branches."
$ws.Range("B10").Value = 35.421093
$ws.Range("C10").Value = 3835
$ws.Range("D10").Value = "0.00789"
$ws.Range("E10").Value = "eba4b5c5-42ac-47c7-afa7-402fe2915a66"

# Row 11
$ws.Range("A11").Value = "Move to location (6, 6) and remove the screws."
$ws.Range("B11").Value = 24.927466
$ws.Range("C11").Value = 3776
$ws.Range("D11").Value = "0.00738"
$ws.Range("E11").Value = "265975ff-7806-4569-8fc0-0dd909dc8afe"

Write-Host "Edits applied"
